# "Upload Module done for excel"
# The MCQ / MTF / FIB answer-option text values that used to be plain text
# are replaced with uploaded image filenames (1.png .. 4.png).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MCQ options for "Which keyword is used for exception handling?")
# super/this/finalize/try -> 1.png/2.png/3.png/4.png
$ws.Range("E2").Value = "1.png"
$ws.Range("F2").Value = "2.png"
$ws.Range("G2").Value = "3.png"
$ws.Range("M2").Value = "4.png"

# Row 3 (MTF "Match the correct pair"): char -> 1.png
$ws.Range("E3").Value = "1.png"

# Row 4 (FIB "Fill In the blanks with proper article"): sentence -> 1.png
$ws.Range("D4").Value = "1.png"

# Move / persist the active selection onto E3, matching the edited cell.
$ws.Range("E3").Select()
